# "export usage et occupations"
#
# Slide 3 ("3 - USAGE ET OCCUPATION DU BATIMENT") carries a single
# free-text "ZoneTexte 2" placeholder. Turn it into the first of three
# stacked, named placeholders (nomCalendrier / nomZones / remarque...)
# that the exporter fills in afterwards, by cloning it twice and
# restacking the trio down the right-hand column of the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# PowerPoint's COM surface measures Left/Top/Width/Height in points
# (1 pt = 12700 EMU) even though the OOXML stores English Metric Units.
# Add a half-EMU nudge before the division so the float round-trip lands
# back on the exact integer EMU value instead of quietly truncating down.
function EMU($emu) {
    return ($emu + 0.5) / 12700.0
}

# Locate the existing "ZoneTexte 2" shape (id=3) on the slide.
$descr = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "ZoneTexte 2") {
        $descr = $cand
    }
}

# 1) Repurpose it as "nomCalendrier" and move/resize it to the top of
#    the right-hand stack; its free text is cleared to a single space.
$descr.Name = "nomCalendrier"
$descr.Left = EMU(6095999)
$descr.Top = EMU(1061595)
$descr.Width = EMU(5810055)
$descr.Height = EMU(338554)
$descr.TextFrame.TextRange.Text = " "

# 2) Clone it for "nomZones", directly below.
$zones = $descr.Duplicate()
$zones.Name = "nomZones"
$zones.Left = EMU(6095999)
$zones.Top = EMU(1622572)
$zones.Width = EMU(5810055)
$zones.Height = EMU(338554)
$zones.TextFrame.TextRange.Text = " "

# 3) Clone it again for the occupation/regulation remark field.
$remarque = $descr.Duplicate()
$remarque.Name = "remarqueRemarques occupation et régulation"
$remarque.Left = EMU(6095998)
$remarque.Top = EMU(2209549)
$remarque.Width = EMU(5810055)
$remarque.Height = EMU(338554)
$remarque.TextFrame.TextRange.Text = " "
